$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = "`n____________<TK>____________ (@tlgkyck) " + [char]0x2022 + " Instagram photos and videos`n"

$ws.Range("A1").Value = $text
$ws.Range("P1").Value = $text
$ws.Rows.Item(1).AutoFit()
